# Add a new "Save" column (H) to the s_vals sheet, mirroring the header
# style used by the existing columns (B1:G1) and filling the data rows
# with 0, matching the other numeric columns already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: same formatting as the other header cells (bold, bordered,
# centered) by copying the format from the neighboring "sum" header (G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Data rows: plain numeric zeros, no special formatting (same as F2:F5/G2:G5).
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
